$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted at row 52 ("Fecha" 2022-11-30 / serial
# 44895). Every existing record from row 52 down to row 140 shifts down by one
# row (row 140 -> row 141), and the brand new record's values are written into
# row 52. We shift bottom-up so we never overwrite data before it's copied,
# and we use Range.Copy so number formats / styles (the date style on column D)
# move along with the values.

for ($r = 140; $r -ge 52; $r--) {
    $src = $ws.Range("A" + $r + ":R" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":R" + ($r + 1))
    $src.Copy($dst)
}

# Write the new record's values into row 52.
$ws.Range("D52").Value2 = 44895
$ws.Range("J52").Value2 = 1100
$ws.Range("K52").Value2 = 21000
$ws.Range("L52").Value2 = 22000
$ws.Range("M52").Value2 = 21500
$ws.Range("P52").Value2 = 307
